$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The Price column (D) holds values that look numeric (and some that use a
# "." thousands separator, e.g. "26.980.18") but the sheet stores all of
# them as plain text. Force the column to Text format first so assigning a
# numeric-looking string below doesn't get silently coerced into a real
# number (which would also normalize away values like trailing zeros).
$ws.Range("D2:D51").NumberFormat = "@"

# Row updates: row number => column letter => new value
$updates = @(
    ,@(2, 4, "26.920.06")
    ,@(2, 5, "  +1.75%  ")
    ,@(3, 4, "1.837.75")
    ,@(3, 5, "  +1.49%  ")
    ,@(4, 4, "1.009")
    ,@(4, 5, "  +0.49%  ")
    ,@(5, 4, "308.82")
    ,@(5, 5, "  +0.93%  ")
    ,@(6, 4, "1.007")
    ,@(6, 5, "  +0.32%  ")
    ,@(7, 4, "0.4656")
    ,@(7, 5, "  +3.40%  ")
    ,@(8, 4, "0.3616")
    ,@(8, 5, "  +0.91%  ")
    ,@(9, 4, "0.07116")
    ,@(9, 5, "  +0.76%  ")
    ,@(10, 4, "0.9093")
    ,@(10, 5, "  +2.34%  ")
    ,@(11, 4, "19.48")
    ,@(11, 5, "  +0.69%  ")
    ,@(12, 4, "0.07670")
    ,@(12, 5, "  -1.45%  ")
    ,@(13, 4, "1.825.45")
    ,@(13, 5, "  +1.11%  ")
    ,@(14, 4, "5.253")
    ,@(14, 5, "  -0.29%  ")
    ,@(15, 4, "6.371")
    ,@(15, 5, "  +0.99%  ")
    ,@(16, 4, "87.91")
    ,@(16, 5, "  +3.72%  ")
    ,@(17, 4, "1.010")
    ,@(18, 4, "0.000008567")
    ,@(18, 5, "  +0.67%  ")
    ,@(19, 4, "1.007")
    ,@(19, 5, "  +0.36%  ")
    ,@(20, 4, "26.967.52")
    ,@(20, 5, "  +1.81%  ")
    ,@(21, 4, "14.26")
    ,@(21, 5, "  +0.65%  ")
    ,@(22, 4, "4.999")
    ,@(22, 5, "  +0.78%  ")
    ,@(23, 4, "10.61")
    ,@(23, 5, "  +0.86%  ")
    ,@(24, 4, "1.926")
    ,@(24, 5, "  -1.30%  ")
    ,@(25, 4, "152.16")
    ,@(25, 5, "  +0.64%  ")
    ,@(26, 4, "18.14")
    ,@(26, 5, "  +2.03%  ")
    ,@(27, 4, "2.023")
    ,@(27, 5, "  -1.73%  ")
    ,@(28, 4, "113.80")
    ,@(28, 5, "  +1.44%  ")
    ,@(29, 4, "4.867")
    ,@(29, 5, "  +0.57%  ")
    ,@(30, 4, "0.08847")
    ,@(30, 5, "  +1.88%  ")
    ,@(31, 4, "3.200")
    ,@(31, 5, "  +2.56%  ")
    ,@(32, 4, "2.801")
    ,@(32, 5, "  +2.16%  ")
    ,@(33, 4, "0.7440")
    ,@(33, 5, "  +0.23%  ")
    ,@(34, 4, "1.165")
    ,@(34, 5, "  +5.06%  ")
    ,@(35, 4, "4.449")
    ,@(35, 5, "  +0.27%  ")
    ,@(36, 4, "1.080")
    ,@(36, 5, "  +0.96%  ")
    ,@(37, 4, "2.972")
    ,@(37, 5, "  +2.82%  ")
    ,@(38, 4, "0.01931")
    ,@(38, 5, "  +0.32%  ")
    ,@(39, 4, "0.05151")
    ,@(39, 5, "  +0.68%  ")
    ,@(40, 4, "0.5149")
    ,@(40, 5, "  +1.38%  ")
    ,@(41, 4, "6.874")
    ,@(41, 5, "  +1.78%  ")
    ,@(42, 4, "0.1508")
    ,@(42, 5, "  +0.05%  ")
    ,@(43, 4, "8.088")
    ,@(43, 5, "  +0.52%  ")
    ,@(44, 4, "10.46")
    ,@(44, 5, "  +4.62%  ")
    ,@(45, 2, "PaxDollar")
    ,@(45, 3, "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp")
    ,@(45, 4, "1.007")
    ,@(45, 5, "  +0.30%  ")
    ,@(46, 2, "Decentraland")
    ,@(46, 3, "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana")
    ,@(46, 4, "0.4665")
    ,@(46, 5, "  -0.28%  ")
    ,@(47, 4, "100.25")
    ,@(47, 5, "  +0.51%  ")
    ,@(48, 4, "1.599")
    ,@(48, 5, "  +1.57%  ")
    ,@(49, 4, "0.06036")
    ,@(49, 5, "  +0.79%  ")
    ,@(50, 4, "64.18")
    ,@(50, 5, "  +0.86%  ")
    ,@(51, 4, "36.09")
    ,@(51, 5, "  +0.59%  ")
)

foreach ($u in $updates) {
    $ws.Cells.Item($u[0], $u[1]).Value = $u[2]
}